$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the candidate's administrative info (row 4) -----------------
# Porteur du projet (A:J)
$ws.Range("A4").Value = "M"
$ws.Range("B4").Value = "Caillau"
$ws.Range("C4").Value = "Jean-Baptiste"
$ws.Range("D4").Value = "jean-baptiste.caillau@univ-cotedazur.fr"
$ws.Range("E4").Value = "Université Côte d’Azur"
$ws.Range("F4").Value = "UMR 7351"
$ws.Range("G4").Value = "LJAD"
$ws.Range("H4").Value = "INSMI"
$ws.Range("J4").Value = "thierry.goudon@univ-cotedazur.fr"

# Partenaire 1 du projet (K:S)
$ws.Range("K4").Value = "M"
$ws.Range("L4").Value = "Sugny"
$ws.Range("M4").Value = "Dominique"
$ws.Range("N4").Value = "dominique.sugny@u-bourgogne.fr"
$ws.Range("O4").Value = "Univeristé Bourgogne Europe"
$ws.Range("P4").Value = "UMR 6303"
$ws.Range("Q4").Value = "ICB"
$ws.Range("R4").Value = "INP"

# Informations projet (AC:AD)
$ws.Range("AC4").Value = "Contrôle théorique et expérimental des centres NV"
$ws.Range("AD4").Value = "CONV"

# Délégations régionales (filled last, like a dropdown pick)
$ws.Range("S4").Value = "Centre Est"
$ws.Range("I4").Value = "Côte d'Azur"

# Budget (AE:AF) with currency number format
$ws.Range("AE4").Value = 11000
$ws.Range("AE4").NumberFormat = '#,##0\ [$€-1];[Red]\-#,##0\ [$€-1]'
$ws.Range("AF4").Value = 10000
$ws.Range("AF4").NumberFormat = '#,##0\ [$€-1];[Red]\-#,##0\ [$€-1]'

# --- Row heights (autofit-driven by the new content) ---------------------
$ws.Rows(3).RowHeight = 69
$ws.Rows(4).RowHeight = 15.75

# --- Selection / view state -----------------------------------------------
$ws.Range("S4").Select()
